# Repull data: update the dSF column (column F) values for specific rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -4
    8  = 0
    9  = -2
    15 = -1
    24 = -8
    29 = 4
    33 = 0
    41 = -4
    42 = -7
    47 = 1
    48 = -1
    51 = 4
    53 = 0
    55 = -1
    73 = -1
    74 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
